$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 502, shifting existing rows 502:549 down to 503:550.
$ws.Rows.Item(502).Insert()

# Populate the newly inserted row 502 with the new weekly price entry.
# Most descriptive columns mirror the row that used to be at 502 (now 503),
# while D (Fecha), J (Volumen), K/L/M (precios) and P (Precio $/Kg) are new.
$ws.Range("A502").Value = 6
$ws.Range("B502").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C502").Value = "Metropolitana"
$ws.Range("D502").Value = 44769
$ws.Range("E502").Value = 13
$ws.Range("F502").Value = 100112044
$ws.Range("G502").Value = "Perejil"
$ws.Range("H502").Value = "Sin especificar"
$ws.Range("I502").Value = "Primera"
$ws.Range("J502").Value = 90
$ws.Range("K502").Value = 19000
$ws.Range("L502").Value = 20000
$ws.Range("M502").Value = 19667
$ws.Range("N502").Value = "$/docena de atados"
$ws.Range("O502").Value = "Región Metropolitana"
$ws.Range("P502").Value = 6556
$ws.Range("Q502").Value = 3
$ws.Range("R502").Value = "Hortaliza"
